$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Modify SLG building config (StateFunction sheet): clear a number of the
# EFT_* flag columns (set to 0) for rows 2-13, keeping ID (col A), EFT_INFO
# (col B) and EFT_FINISH (col O) as-is, along with a few per-row exceptions
# that remain 1.
$ws.Range("E2:N2").Value = 0
$ws.Range("C3:D3").Value = 0
$ws.Range("F3:N3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4:N4").Value = 0
$ws.Range("C5:D5").Value = 0
$ws.Range("F5:N5").Value = 0
$ws.Range("C6:N6").Value = 0
$ws.Range("C7:N7").Value = 0
$ws.Range("C8:N8").Value = 0
$ws.Range("C9:N9").Value = 0
$ws.Range("C10:N10").Value = 0
$ws.Range("C11:N11").Value = 0
$ws.Range("C12:N12").Value = 0
$ws.Range("C13:N13").Value = 0

# Active selection ends up on F11 in the saved file.
$ws.Range("F11").Select()
